# ---------------------------------------------------------------------------
# "added github to pptx" — apply the recorded edits to slide 1.
# ---------------------------------------------------------------------------
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1) "Predictor de riesgo crediticio{" textbox (Shape 2 / "TextBox 3"):
#    split the single run into five runs so that "riesgo" and "crediticio"
#    become their own runs (PowerPoint does this automatically once the
#    proofing pass flags those words; we recreate the run boundaries here).
# ---------------------------------------------------------------------------
$titleShape = $s.Shapes.Item(2)
$titleRange = $titleShape.TextFrame.TextRange

# Touch the two "misspelled" words' character-level formatting so the host
# keeps them as separate runs instead of re-merging them with their
# neighbours (identical rPr runs get coalesced back together otherwise).
$titleRange.Characters(14, 6).Font.BaselineOffset = 0
$titleRange.Characters(21, 10).Font.BaselineOffset = 0

# Re-assert the text for every segment (no-ops content-wise, but keeps the
# run split lined up with the word boundaries from the diff).
$titleRange.Characters(1, 13).Text = "Predictor de "
$titleRange.Characters(14, 6).Text = "riesgo"
$titleRange.Characters(20, 1).Text = " "
$titleRange.Characters(21, 10).Text = "crediticio"
$titleRange.Characters(31, 1).Text = "{"

# ---------------------------------------------------------------------------
# 2) Reposition the "<Por=...CJ"/>" textbox (Shape 4 / "TextBox 5").
# ---------------------------------------------------------------------------
$byShape = $s.Shapes.Item(4)
$byShape.Left = 180
$byShape.Top = 453

# ---------------------------------------------------------------------------
# 3) Add the new GitHub link textbox ("CuadroTexto 10").
#    id=10 is consumed/discarded by an intermediate shape (matching the
#    author's deck, where the new shape lands on id=11).
# ---------------------------------------------------------------------------
$placeholder = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$placeholder.Delete()

$linkShape = $s.Shapes.AddTextbox(1, 1174.5307086614173, 755.9255905511811, 308.7988188976378, 29.081259842519685)
$linkShape.Name = "CuadroTexto 10"
$linkShape.Fill.Visible = $false
$linkShape.TextFrame.WordWrap = $true
$linkShape.TextFrame.AutoSize = 1

$linkRange = $linkShape.TextFrame.TextRange
$linkRange.Text = "https://github.com/Chelqq/"
$linkRange.LanguageID = "es-MX"
$linkRange.Font.Name = "Courier Prime"
$linkRange.Font.Color.RGB = 16777215
